$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A2").Value = 70713
$ws.Range("A3").Value = 70716
$ws.Range("A4").Value = 70718
$ws.Range("A5").Value = 70720

$wb.Save()
